$wb = $excel.ActiveWorkbook

# "invalidcode" is the last existing sheet and has the right base sheet
# format (no <cols>, defaultRowHeight 14.5, tabSelected). Duplicate it so
# the new sheet inherits that formatting, then rename/re-populate it.
$src = $wb.Worksheets.Item("invalidcode")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy([System.Reflection.Missing]::Value, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "tryeditorcode"

# Pull matching cell formatting (same style combos already used on the
# "validcode" sheet) over to the new cells before filling in values, so
# the saved workbook reuses the existing style indices instead of minting
# new ones.
$vc = $wb.Worksheets.Item("validcode")

$vc.Range("B1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)

$vc.Range("B2").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)

$vc.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "pythonCode"
$newSheet.Range("B1").Value = "output"
$newSheet.Range("A2").Value = "print(""Hello"")"
$newSheet.Range("B2").Value = "Hello"
$newSheet.Range("A3").Value = "Hello"

$newSheet.Range("A1:B4").Select() | Out-Null
